$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the subtitle text to reflect the new Year-to-Date month (October -> November 2016)
$ws.Range("A2").Value = "Total (All Sectors) by Census Division and State, Year-to-Date through November 2016 (Continued)"

# Updated Relative Standard Error data cells (EPM_2016_11 monthly refresh)
$ws.Range("E4").Value = 11
$ws.Range("H4").Value = 5
$ws.Range("E5").Value = 55
$ws.Range("F5").Value = 5
$ws.Range("H5").Value = 7
$ws.Range("F6").Value = 1
$ws.Range("H6").Value = 14
$ws.Range("I6").Value = 5
$ws.Range("E7").Value = 11
$ws.Range("H7").Value = 7
$ws.Range("I7").Value = 3
$ws.Range("F8").Value = 11
$ws.Range("H8").Value = 45
$ws.Range("E9").Value = 62
$ws.Range("F9").Value = 17
$ws.Range("E10").Value = 34
$ws.Range("I10").Value = 22
$ws.Range("E11").Value = 9
$ws.Range("H11").Value = 5
$ws.Range("E12").Value = 10
$ws.Range("F12").Value = 7
$ws.Range("H12").Value = 8
$ws.Range("H13").Value = 8
$ws.Range("E14").Value = 30
$ws.Range("E15").Value = 13
$ws.Range("I15").Value = 0.31
$ws.Range("E16").Value = 35
$ws.Range("H16").Value = 27
$ws.Range("E17").Value = 16
$ws.Range("E18").Value = 58
$ws.Range("H18").Value = 14
$ws.Range("E19").Value = 31
$ws.Range("F19").Value = 4
$ws.Range("H19").Value = 42
$ws.Range("E20").Value = 201
$ws.Range("F20").Value = 3
$ws.Range("H20").Value = 41
$ws.Range("E21").Value = 40
$ws.Range("F21").Value = 0.36
$ws.Range("H21").Value = 13
$ws.Range("F22").Value = 0.42
$ws.Range("H22").Value = 187
$ws.Range("E23").Value = 204
$ws.Range("F23").Value = 0.29
$ws.Range("E24").Value = 101
$ws.Range("H24").Value = 13
$ws.Range("E25").Value = 48
$ws.Range("E26").Value = 102
$ws.Range("F27").Value = 0.37
$ws.Range("H27").Value = 54
$ws.Range("I27").Value = 2
$ws.Range("I28").Value = 1
$ws.Range("I29").Value = 0.22
$ws.Range("E30").Value = 35
$ws.Range("I30").Value = 6
$ws.Range("I31").Value = 156
$ws.Range("I32").Value = 1
$ws.Range("F33").Value = 2
$ws.Range("I33").Value = 0.39
$ws.Range("E34").Value = 19
$ws.Range("E35").Value = 5
$ws.Range("H35").Value = 16
$ws.Range("I35").Value = 0.47
$ws.Range("E36").Value = 121
$ws.Range("H36").Value = 22
$ws.Range("I36").Value = 0.45
$ws.Range("E37").Value = 6
$ws.Range("H37").Value = 6
$ws.Range("I37").Value = 0.46
$ws.Range("F38").Value = 0.47
$ws.Range("I38").Value = 0.49
$ws.Range("E39").Value = 8
$ws.Range("F39").Value = 2
$ws.Range("H39").Value = 14
$ws.Range("I39").Value = 0.4
$ws.Range("F40").Value = 3
$ws.Range("F41").Value = 7
$ws.Range("H42").Value = 179
$ws.Range("E43").Value = 29
$ws.Range("F43").Value = 6
$ws.Range("E44").Value = 6
$ws.Range("F44").Value = 0.35
$ws.Range("H44").Value = 12
$ws.Range("I44").Value = 0.29
$ws.Range("E45").Value = 59
$ws.Range("F45").Value = 3
$ws.Range("H46").Value = 15
$ws.Range("F47").Value = 0.31
$ws.Range("H47").Value = 47
$ws.Range("E48").Value = 6
$ws.Range("F48").Value = 0.44
$ws.Range("H48").Value = 14
$ws.Range("I48").Value = 0.38
$ws.Range("H49").Value = 7
$ws.Range("I49").Value = 1
$ws.Range("I50").Value = 0.18
$ws.Range("E51").Value = 5
$ws.Range("F51").Value = 0.43
$ws.Range("H51").Value = 59
$ws.Range("C52").Value = 48
$ws.Range("F52").Value = 4
$ws.Range("H52").Value = 49
$ws.Range("I52").Value = 6
$ws.Range("F53").Value = 2
$ws.Range("F54").Value = 4
$ws.Range("C55").Value = 122
$ws.Range("E55").Value = 6
$ws.Range("H55").Value = 2158
$ws.Range("C56").Value = 9
$ws.Range("F56").Value = 3
$ws.Range("I56").Value = 2
$ws.Range("E60").Value = 24
$ws.Range("H60").Value = 49
$ws.Range("H61").Value = 43
$ws.Range("E62").Value = 28
$ws.Range("F63").Value = 22
$ws.Range("E64").Value = 28
$ws.Range("I65").Value = 0.16
